$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")

# Revert the password for row 7 (uid=5, nric=S1234567A) from "NEWPASSWORD" back to "Password"
$ws.Range("C7").Value = "Password"
